$d = $word.ActiveDocument

# Update the date heading (unique text, safe to use Find/Replace)
$d.Content.Find.Execute("2025-04-03 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-04 Friday", 2)

# Update the multiplication problems in the table.
# Cells are addressed directly by (row, column) to avoid any ambiguity
# caused by a new value for one cell matching the old value of another.
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text  = "43×34=1462"
$t.Cell(1,2).Range.Text  = "28×50=1400"
$t.Cell(1,3).Range.Text  = "81×14=1134"
$t.Cell(1,4).Range.Text  = "74×88=6512"
$t.Cell(1,5).Range.Text  = "31×85=2635"

$t.Cell(5,1).Range.Text  = "29×95=2755"
$t.Cell(5,2).Range.Text  = "15×34=510"
$t.Cell(5,3).Range.Text  = "76×97=7372"
$t.Cell(5,4).Range.Text  = "55×71=3905"
$t.Cell(5,5).Range.Text  = "71×33=2343"

$t.Cell(10,1).Range.Text = "97×62=6014"
$t.Cell(10,2).Range.Text = "12×55=660"
$t.Cell(10,3).Range.Text = "52×65=3380"
$t.Cell(10,4).Range.Text = "11×83=913"
$t.Cell(10,5).Range.Text = "80×68=5440"

$t.Cell(15,1).Range.Text = "72×55=3960"
$t.Cell(15,2).Range.Text = "52×93=4836"
$t.Cell(15,3).Range.Text = "53×26=1378"
$t.Cell(15,4).Range.Text = "70×16=1120"
$t.Cell(15,5).Range.Text = "11×29=319"

$t.Cell(20,1).Range.Text = "76×40=3040"
$t.Cell(20,2).Range.Text = "29×68=1972"
$t.Cell(20,3).Range.Text = "33×77=2541"
$t.Cell(20,4).Range.Text = "80×74=5920"
$t.Cell(20,5).Range.Text = "66×41=2706"
